# Update the lattice multiplication table with new problems/digits.
# Each table cell holds one lattice-multiplication problem laid out as
# five text runs separated by manual line breaks (char 11 / <w:br/>):
#   "N1 x N2"
#   "  d1    d2"   (digits of N2)
#   "  ----"
#   "d3|    |"     (tens digit of N1)
#   "d4|    |"     (ones digit of N1)

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$cell = $table.Cell(1, 1)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>41 x 53</w:t><w:br/><w:t xml:space="preserve">  5    3</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>4|    |</w:t><w:br/><w:t>1|    |</w:t></w:r></w:p>')

$cell = $table.Cell(1, 2)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>13 x 76</w:t><w:br/><w:t xml:space="preserve">  7    6</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>3|    |</w:t></w:r></w:p>')

$cell = $table.Cell(1, 3)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>17 x 39</w:t><w:br/><w:t xml:space="preserve">  3    9</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>7|    |</w:t></w:r></w:p>')

$cell = $table.Cell(2, 1)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>57 x 30</w:t><w:br/><w:t xml:space="preserve">  3    0</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>5|    |</w:t><w:br/><w:t>7|    |</w:t></w:r></w:p>')

$cell = $table.Cell(2, 2)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>12 x 90</w:t><w:br/><w:t xml:space="preserve">  9    0</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>2|    |</w:t></w:r></w:p>')

$cell = $table.Cell(2, 3)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>67 x 20</w:t><w:br/><w:t xml:space="preserve">  2    0</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>7|    |</w:t></w:r></w:p>')

$cell = $table.Cell(3, 1)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>70 x 87</w:t><w:br/><w:t xml:space="preserve">  8    7</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>7|    |</w:t><w:br/><w:t>0|    |</w:t></w:r></w:p>')

$cell = $table.Cell(3, 2)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>14 x 43</w:t><w:br/><w:t xml:space="preserve">  4    3</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>4|    |</w:t></w:r></w:p>')

$cell = $table.Cell(3, 3)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>96 x 55</w:t><w:br/><w:t xml:space="preserve">  5    5</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>9|    |</w:t><w:br/><w:t>6|    |</w:t></w:r></w:p>')

$cell = $table.Cell(4, 1)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>92 x 34</w:t><w:br/><w:t xml:space="preserve">  3    4</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>9|    |</w:t><w:br/><w:t>2|    |</w:t></w:r></w:p>')

$cell = $table.Cell(4, 2)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>73 x 20</w:t><w:br/><w:t xml:space="preserve">  2    0</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>7|    |</w:t><w:br/><w:t>3|    |</w:t></w:r></w:p>')

$cell = $table.Cell(4, 3)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>15 x 88</w:t><w:br/><w:t xml:space="preserve">  8    8</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p>')

$cell = $table.Cell(5, 1)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>95 x 83</w:t><w:br/><w:t xml:space="preserve">  8    3</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>9|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p>')

$cell = $table.Cell(5, 2)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>35 x 32</w:t><w:br/><w:t xml:space="preserve">  3    2</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>3|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p>')

$cell = $table.Cell(5, 3)
$cell.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>91 x 92</w:t><w:br/><w:t xml:space="preserve">  9    2</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>9|    |</w:t><w:br/><w:t>1|    |</w:t></w:r></w:p>')
